$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Part 1: "Alpha Beta Pruning is commonly used in two player games a|nd
# was considered..." paragraph.
#
# In the starting document this sentence is split into two runs around the
# (old) "_GoBack" bookmark. Running Find/Replace across that split rejoins
# the text into a single run and removes the bookmark that lived there,
# matching the target XML.
# -----------------------------------------------------------------------
$rAlpha = $d.Content
$rAlpha.Find.ClearFormatting()
$rAlpha.Find.Execute("games and was considered", $false, $false, $false, $false, $false, $true, 1, $false, "games and was considered", 2)

# -----------------------------------------------------------------------
# Part 2: "This strategy can result in the AI taking fewer hits..."
# paragraph.
#
# "can" becomes "should", and the "_GoBack" bookmark moves here, right
# after the freshly-typed word - mirroring where Word leaves the cursor
# bookmark after an edit. The target XML keeps this run split into three
# pieces ("This strategy " / "should" / " result in the AI taking ") even
# though all three share identical formatting, so we plant short-lived
# "splitter" bookmarks at each boundary *before* mutating the text. Any
# text mutation triggers a run-consolidation pass across the whole
# paragraph, but a bookmark sitting exactly on a boundary keeps the runs
# on each side of it distinct even after that temporary bookmark is
# deleted again afterwards.
# -----------------------------------------------------------------------
$rThis = $d.Content
$rThis.Find.ClearFormatting()
$rThis.Find.Execute("This strategy can result in the AI taking ")
$sentenceStart = $rThis.Start
$sentenceEnd = $rThis.End

$boundaryBeforeSentence = $d.Range($sentenceStart, $sentenceStart)
$d.Bookmarks.Add("ZZZTempBoundaryStart", $boundaryBeforeSentence)

$canRange = $d.Range($sentenceStart + 14, $sentenceStart + 17)
$boundaryBeforeCan = $d.Range($canRange.Start, $canRange.Start)
$d.Bookmarks.Add("ZZZTempBoundaryMid", $boundaryBeforeCan)

$boundaryAfterSentence = $d.Range($sentenceEnd, $sentenceEnd)
$d.Bookmarks.Add("ZZZTempBoundaryEnd", $boundaryAfterSentence)

$canRange.Text = "should"

$afterShould = $d.Range($canRange.End, $canRange.End)
$d.Bookmarks.Add("_GoBack", $afterShould)

$d.Bookmarks("ZZZTempBoundaryStart").Delete()
$d.Bookmarks("ZZZTempBoundaryMid").Delete()
$d.Bookmarks("ZZZTempBoundaryEnd").Delete()
